$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (typo / wording fixes) ---
$ws.Range("B10").Value = "Are we allowed to change the source code by our own, if we want to ?"
$ws.Range("C12").Value = "How ever you want. There is no restriction, you will get contact data of all involved persons."
$ws.Range("B15").Value = "How much time do you need in advance if we want to give further jobs ?"

# --- Formatting: give the question/answer area a white background fill ---
$ws.Range("B2:G22").Interior.Color = 16777215

# --- Remove the frozen header row/column pane ---
$excel.ActiveWindow.FreezePanes = $false
